$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (D1/E1 renamed, G1 added) ---
$ws.Range("D1").Value = "Debit"
$ws.Range("E1").Value = "Credit"

# G1 is a brand new header cell; copy the header formatting from F1
# (bold font, thin border, center/top alignment) before giving it a value
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("G1").Value = "Reconciled"

# --- Row 2 ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2/9/25"
$ws.Range("B2").Value = "paycheck"
$ws.Range("C2").Value = "me"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 400
$ws.Range("F2").Value = 400

# --- Row 3 ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2/9/25"
$ws.Range("B3").Value = "Bar Tab"
$ws.Range("C3").Value = "Some Watering Hole"
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 370

# --- Row 4 ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2/9/25"
$ws.Range("B4").Value = "Dinner"
$ws.Range("C4").Value = "Slim Chicken"
$ws.Range("D4").Value = 25.69
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 344.31

# --- Row 5 ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2/9/25"
$ws.Range("B5").Value = "Chime Transfer"
$ws.Range("C5").Value = "Chime"
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 244.31
$ws.Range("G5").Value = 1

# --- Row 6 ---
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-03-01"
$ws.Range("B6").Value = "Starting Balance"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 244.31
$ws.Range("G6").Value = 0

# --- Row 7 ---
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2/10/25"
$ws.Range("B7").Value = "Some Description"
$ws.Range("C7").Value = "A Payee"
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 236.31
$ws.Range("G7").Value = 1

# --- Row 8 ---
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "3/10/25"
$ws.Range("B8").Value = "Some Description"
$ws.Range("C8").Value = "A Payee"
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 228.31
$ws.Range("G8").Value = 1
